# test-23: extrapol 2.0 & hybrid3 2.0 (mae)
# Fill in the MAE "extrapol 2.0" / "Hybrid3 2.0" results on the
# "mae (3Ysum)" sheet, and clear out the placeholder numbers that had
# been left (mistakenly, under the non-"2.0" labels) on "mse (3Ysum)",
# correcting that sheet's headers to the "2.0" wording too.

$wb = $excel.ActiveWorkbook

$wsMae = $wb.Worksheets.Item("mae (3Ysum)")
$wsMse = $wb.Worksheets.Item("mse (3Ysum)")

# ---- mse (3Ysum): fix headers to the "2.0" variants, clear its K/P data ----
$wsMse.Range("J3").Value = "RF-100 (superdataset-24-f 2Y.csv + extrapol 2.0)"
$wsMse.Range("O3").Value = "Hybrid3-model 2.0 (superdataset-24-f + 2Y + 3Y.csv)"
$wsMse.Range("K5:K54").ClearContents()
$wsMse.Range("P5:P54").ClearContents()

# ---- mae (3Ysum): fill in the new MAE results for both models ----
$maeK = @(206.66380660954709,212.8558629130967,214.32481028151781,201.02997552019579,197.1347735618115,195.36096695226439,204.97538555691551,204.89599755201959,195.98128518971851,200.9626438188495,213.92667074663399,206.2847246022032,206.01085679314559,217.98350061199511,205.23118727050181,187.24915544675639,196.88682986536111,214.93334149326799,205.9227539779682,198.87408812729501,197.95283965728271,207.66452876376991,197.90394124847001,209.68507955936349,202.95128518971839,219.02046511627901,197.50045287637701,194.2596205630355,207.83381884944919,205.1839167686658,199.70571603427169,199.8348347613219,201.90072215422271,210.15872705018359,197.1021787025704,208.03127294981641,210.57657282741741,202.71604651162789,198.43517747858019,203.32753977968181,212.9602447980416,203.55583843329251,204.4150305997552,207.16824969400241,195.307429620563,198.4937576499388,210.52051407588741,195.04549571603431,206.35179926560579,200.96567931456551)
$maeP = @(201.90884944920441,208.26321909424721,202.3569522643819,195.578029375765,204.21499388004901,193.42537331701351,194.7519094247246,209.19140758873931,189.51013463892289,196.13413708690331,199.1854957160343,195.36242350061201,199.86490820073439,202.64277845777229,189.8554345165239,197.23663402692779,203.65611995104041,198.8673561811506,201.3649938800489,189.76112607099139,200.21166462668299,193.31332925336599,186.67208078335369,201.99212974296211,186.95958384332931,193.643341493268,197.10949816401461,203.89558139534881,200.75492044063651,182.58177478580171,193.973476132191,192.87325581395351,201.4863892288862,189.74165238678091,198.32725826193391,191.91577723378211,202.66350061199509,201.27548347613219,202.94133414932679,211.36705018359851,209.0415544675642,198.20809057527541,196.59665850673201,208.57662178702569,202.00973072215419,192.64286413708689,190.80787025703791,197.01254589963281,200.15489596083231,194.21965728274171)

for ($i = 0; $i -lt $maeK.Length; $i++) {
    $row = 5 + $i
    $wsMae.Cells.Item($row, 11).Value = $maeK[$i]
    $wsMae.Cells.Item($row, 16).Value = $maeP[$i]
}

# ---- selection / active-tab bookkeeping ----
$wsMae.Range("S39").Select()
$wsMse.Select()
$wsMse.Range("P18").Select()
